$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the retrospective answers gathered during the meeting
# (order matches how the shared-string table was originally populated)
$ws.Range("A2").Value = "Github anduvo bien toda la semana hasta el ultimo momento"
$ws.Range("C3").Value = "primero hacer commit y dspues sincronizar"
$ws.Range("C4").Value = "la idea es no tocar lo mismo"
$ws.Range("C2").Value = "mejorar github"

# C2's border formatting was cleared when the text was entered
$ws.Range("C2").Borders.LineStyle = -4142

# Leave the selection on C2, matching where the user ended up
$ws.Range("C2").Select()
